$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as sheet row 34 (just above the
# record that used to be row 34), pushing every subsequent record down by one
# row. Insert a blank row at 34 first so everything below shifts down, then
# populate the new row with the new observation's data.
$ws.Rows("34:34").Insert()

$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44526
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = 100112028
$ws.Range("G34").Value = "Sandia"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 600
$ws.Range("K34").Value = 1300
$ws.Range("L34").Value = 1300
$ws.Range("M34").Value = 1300
$ws.Range("N34").Value = '$/kilo (volumen en unidades)'
$ws.Range("O34").Value = "Perú"
$ws.Range("P34").Value = 1300
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"
